$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the mailto hyperlinks that live on the row-7 record before the
# row itself is deleted, so the relationship/hyperlink entries for
# bittonnir12@gmail.com / nevilgreen@gmail.com don't linger.
$ws.Range("C7").Hyperlinks.Delete()
$ws.Range("D7").Hyperlinks.Delete()

# Delete the entire row 7 (the review-database record for
# bittonnir12@gmail.com / nevilgreen@gmail.com). This shifts the
# trailing blank template row (old row 8) up to become the new row 7,
# matching the target sheet's A1:F7 dimension.
$ws.Rows.Item(7).Delete()

# Move the active selection to A7, matching the post-edit cursor
# position left behind in the workbook.
$ws.Range("A7").Select()
